# Update the COVID-19 Valais daily figures for rows 469-472
# (new case counts for 2021-06-08 .. 2021-06-10, and fill in 2021-06-11).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 469 (2021-06-08): "Nb nouveaux cas positifs" 14 -> 15
$ws.Range("C469").Value = 15

# --- Row 470 (2021-06-09): "Nb nouveaux cas positifs" 4 -> 18
$ws.Range("C470").Value = 18

# --- Row 471 (2021-06-10): "Nb nouveaux cas positifs" 0 -> 9
$ws.Range("C471").Value = 9

# L471/M471 were storing the text "0" (shared string) - turn them into
# genuine numeric zeros while keeping their existing (text) display format.
$ws.Range("L471").NumberFormat = "General"
$ws.Range("L471").Value = 0
$ws.Range("L471").NumberFormat = "@"

$ws.Range("M471").NumberFormat = "General"
$ws.Range("M471").Value = 0
$ws.Range("M471").NumberFormat = "@"

# --- Row 472 (2021-06-11): fill in the previously-empty daily entry with
# zero new cases (C472 stays blank, L472/M472 become numeric zeros).
$ws.Range("L472").NumberFormat = "General"
$ws.Range("L472").Value = 0
$ws.Range("L472").NumberFormat = "@"

$ws.Range("M472").NumberFormat = "General"
$ws.Range("M472").Value = 0
$ws.Range("M472").NumberFormat = "@"

# --- Sheet view: scroll frozen pane back to the top and reselect A2
$ws.Range("A2").Select() | Out-Null
